{"js": "// Waseef report assignment update\n// 1) Remove the six scoping bullet paragraphs right after \"Talk about the\n//    following:\" (Concepts.../Example.../Data:.../Brief introduction.../\n//    Setting up Geth/Open the Notebook...), leaving \"Talk about the\n//    following:\" immediately followed by \"Talk about timeseries, ...\".\nconst body = context.document.body;\n\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst toDelete = [];\nfor (const p of paras.items) {\n  const t = p.text;\n  if (\n    t === \"Concepts: Ethereum, Tokens, Smart contracts, ERC20\" ||\n    t === \"Example of a smart contract\" ||\n    t ===\n      \"Data: Blockchain is stored in leveldb binary format. No JSON or any text format available\" ||\n    t === \"Brief introduction about Geth. \" ||\n    t === \"Setting up Geth\" ||\n    t === \"Open the Notebook Analysis_of_Cryptocurrency_Investments\"\n  ) {\n    toDelete.push(p);\n  }\n}\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n\n// 2) Merge the trailing two runs of the ADF sentence into one run (same\n//    formatting already), turning \" test\" + \". You will find some info\n//    about it on the notebook. Please expand it a little\" into a single\n//    run \" test. You will find some info about it on the notebook.\n//    Please expand it a little\".\nconst tailResults = body.search(\n  \". You will find some info about it on the notebook. Please expand it a little\",\n  { matchCase: true }\n);\ntailResults.load(\"items\");\nawait context.sync();\ntailResults.items[0].delete();\nawait context.sync();\n\nconst testResults = body.search(\" test\", { matchCase: true });\ntestResults.load(\"items\");\nawait context.sync();\nconst testRange = testResults.items[0];\nconst testEnd = testRange.getRange(\"End\");\ntestEnd.insertText(\n  \". You will find some info about it on the notebook. Please expand it a little\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 3) Drop the stale <w:lastRenderedPageBreak/> on the run that starts\n//    \"Include the code of the ARIMA model right below...\". Force a real\n//    content mutation on that run (Word only clears stale page-break\n//    markers on an actual text edit) by nudging the text and restoring it.\nconst arimaResults = body.search(\n  \"Include the code of the ARIMA model right below. DO NOT explain the output. Still need to work on it\",\n  { matchCase: true }\n);\narimaResults.load(\"items\");\nawait context.sync();\nconst arimaRange = arimaResults.items[0];\nconst arimaStart = arimaRange.getRange(\"Start\");\narimaStart.insertText(\"\\u200b\", \"Replace\");\nawait context.sync();\n\nconst arimaResults2 = body.search(\n  \"\\u200bInclude the code of the ARIMA model right below. DO NOT explain the output. Still need to work on it\",\n  { matchCase: true }\n);\narimaResults2.load(\"items\");\nawait context.sync();\narimaResults2.items[0].insertText(\n  \"Include the code of the ARIMA model right below. DO NOT explain the output. Still need to work on it\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 4) Merge the \"Install the code format add on from MS Word. \" paragraph\n//    with the following (bookmark-only) paragraph, removing the paragraph\n//    break between them so the _GoBack bookmark lives in the same\n//    paragraph as the preceding sentence.\nconst finalParas = body.paragraphs;\nfinalParas.load(\"items/text\");\nawait context.sync();\n\nlet installPara = null;\nlet installIndex = -1;\nconst items = finalParas.items;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"Install the code format add on from MS Word. \") {\n    installPara = items[i];\n    installIndex = i;\n    break;\n  }\n}\nconst bookmarkPara = items[installIndex + 1];\nconst joinRange = installPara.getRange(\"End\").expandTo(bookmarkPara.getRange(\"Start\"));\njoinRange.delete();\nawait context.sync();\n", "ps1": "# Waseef report assignment update\n$d = $word.ActiveDocument\n\n# 1) Remove the six scoping bullet paragraphs right after \"Talk about the\n#    following:\" (Concepts.../Example.../Data:.../Brief introduction.../\n#    Setting up Geth/Open the Notebook...), leaving \"Talk about the\n#    following:\" immediately followed by \"Talk about timeseries, ...\".\n$start = $d.Paragraphs(2).Range.Start\n$end = $d.Paragraphs(7).Range.End\n$d.Range($start, $end).Delete()\n\n# 2) Merge the trailing two runs of the ADF sentence into one run (same\n#    formatting already), turning the non-breaking-space+\"test\" run and\n#    the \". You will find some info about it on the notebook. Please\n#    expand it a little\" run into a single run.\n$tail = $d.Content\n$tail.Find.Execute(\". You will find some info about it on the notebook. Please expand it a little\", $true)\n$tail.Delete()\n\n$nbsp = [char]0x00A0\n$testRange = $d.Content\n$testRange.Find.Execute($nbsp + \"test\", $true)\n$testRange.Collapse(0)\n$testRange.Text = \". You will find some info about it on the notebook. Please expand it a little\"\n\n# 3) Drop the stale lastRenderedPageBreak marker on the run that starts\n#    \"Include the code of the ARIMA model right below...\". Word only\n#    clears a stale page-break marker on an actual text edit, so nudge the\n#    text with a temporary character and then restore it.\n$arimaText = \"Include the code of the ARIMA model right below. DO NOT explain the output. Still need to work on it\"\n$arimaRange = $d.Content\n$arimaRange.Find.Execute($arimaText, $true)\n$arimaRange.Text = $arimaText + \"X\"\n\n$arimaRange2 = $d.Content\n$arimaRange2.Find.Execute($arimaText + \"X\", $true)\n$arimaRange2.Text = $arimaText\n\n# 4) Merge the \"Install the code format add on from MS Word. \" paragraph\n#    with the following (bookmark-only) paragraph, removing the paragraph\n#    mark between them so the _GoBack bookmark lives in the same\n#    paragraph as the preceding sentence.\n$installText = \"Install the code format add on from MS Word. \"\n$installRange = $d.Content\n$installRange.Find.Execute($installText, $true)\n$installRange.Collapse(0)\n$installRange.MoveEnd(1, 1)\n$installRange.Delete()\n"}
